$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct / clean up text in a handful of title cells ---
# (ordered to match the original author's edit order so new shared-string
# entries land in the same append order)

# Row 2: drop the misspelled "vineyears" -> "vineyards"
$ws.Range("B2").Value2 = "Forest land at the top of this California hill has been cleared to make room for wine vineyards"

# Row 70: fix "swarfs" -> "dwarfs"
$ws.Range("B70").Value2 = "A tall tree dwarfs a vineyard in various stages of development"

# Row 92: fix "cAlifornia" -> "California"
$ws.Range("B92").Value2 = "A California vineyard located next to a winery in a new winegrowing area"

# Row 11: drop stray "Max Yavno." credit prefix
$ws.Range("B11").Value2 = "Autumn nights turn the leaves of this California vine almost the color of the table wines produced from its grapes"

# Row 4: drop stray "Max Yavno." credit prefix and fix "Vinyards" -> "Vineyards"
$ws.Range("B4").Value2 = "Vineyards are carved from the wooded slopes and peaks of northern California hills"

# Row 2: catalog number typo fix
$ws.Range("A2").Value2 = "B-000c"

# --- Reset the view: scroll back to the top and select B14 ---
$ws.Range("A1").Select() | Out-Null
$ws.Range("B14").Select() | Out-Null
